# Update the "Förändrad" (Changed) date column C for all data rows (2-33)
# from serial date 46060 (2026-02-07) to 46061 (2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
